# Append/refresh edit for "ランサーズ" (Lancers) sheet — 2025-11-15 06:24 JST run.
# - Rows 2-10 get new scrape results (most rows 11-19 from the previous run are gone).
# - Column widths B/D/H shrink slightly to fit the new content.
# - Hyperlinks on F2:F10 are rebuilt to point at the new job URLs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Drop every existing hyperlink up front; we rebuild the ones we still need
# (F2:F10) after the row positions are final, since row deletion does not
# shift existing hyperlink references in this engine.
$ws.Range("A1").Hyperlinks.Delete()

# --- Row 2 ---
$ws.Range("A2").Value = "2025-11-15 06:24:44"
$ws.Range("B2").Value = "【急募】生成AI・RAG活用の業務ナレッジ検索システム改善"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5434552"
$ws.Range("G2").Value = 318
$ws.Range("H2").Value = "🔥AI,Ai"

# --- Row 3 ---
$ws.Range("A3").Value = "2025-11-15 06:24:44"
$ws.Range("B3").Value = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Range("D3").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5251319"
$ws.Range("G3").Value = 135
$ws.Range("H3").Value = "◆ツール,スクレイピング ◇サイト"

# --- Row 4 ---
$ws.Range("A4").Value = "2025-11-15 06:24:44"
$ws.Range("B4").Value = "【システム開発】FileMaker Proを活用した販売システム構築"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5434428"
$ws.Range("G4").Value = 118
$ws.Range("H4").Value = "◆開発,システム開発"

# --- Row 5 ---
$ws.Range("A5").Value = "2025-11-15 06:24:44"
$ws.Range("B5").Value = "初回 あるサイトの自動操作スクリプト開発(作業見積5時間以内/予算1万以内)の仕事・依頼"
$ws.Range("D5").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5434568"
$ws.Range("G5").Value = 85
$ws.Range("H5").Value = "◆開発 ◇サイト"

# --- Row 6 ---
$ws.Range("A6").Value = "2025-11-15 06:24:44"
$ws.Range("B6").Value = "Flutter iOSアプリにおけるRevenueCat導入のバグ修正依頼"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5434437"
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = "◇アプリ"

# --- Row 7 ---
$ws.Range("A7").Value = "2025-11-15 06:24:44"
$ws.Range("B7").Value = "【急募】料理教室のレシピデジタル化とマイページ構築依頼"
$ws.Range("D7").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5434648"
$ws.Range("G7").Value = 18
$ws.Range("H7").ClearContents()

# --- Row 8 ---
$ws.Range("A8").Value = "2025-11-15 06:24:44"
$ws.Range("B8").Value = "【急募】TradingViewインジシグナルを用いたXAUUSD自動売買EA制作"
$ws.Range("D8").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5434524"
$ws.Range("G8").Value = 18
$ws.Range("H8").ClearContents()

# --- Row 9 ---
$ws.Range("A9").Value = "2025-11-15 06:24:44"
$ws.Range("B9").Value = "URL付きPDF資料の閲覧状況を可視化し、トラッキングする"
$ws.Range("D9").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5434431"
$ws.Range("G9").Value = 18
$ws.Range("H9").ClearContents()

# --- Row 10 ---
$ws.Range("A10").Value = "2025-11-15 06:24:44"
$ws.Range("B10").Value = "GAS構築できる方求む!"
$ws.Range("D10").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5434226"
$ws.Range("G10").Value = 13
$ws.Range("H10").ClearContents()

# Remove the now-stale rows 11-19 (old results that dropped off this run).
$ws.Range("A11:A19").EntireRow.Delete()

# Column width tweaks to fit the refreshed copy.
# (.ColumnWidth is in "characters"; the saved OOXML <col width> is
# ColumnWidth + 5/6, so back that constant off to land on the target widths.)
$ws.Columns.Item(2).ColumnWidth = 51 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 28 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 19 - (5/6)

# Rebuild the hyperlinks for the URL column now that rows 2-10 are final.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5434552")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5251319")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5434428")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5434568")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5434437")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5434648")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5434524")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5434431")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5434226")
